$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D23 text first so it becomes shared string index 19
$ws.Range("D23").Value = "Put in the basic structure for a menu, with File, Edit and Run. Discovered Java doesn’t cope with high dpi displays without manually setting font sizes."
# Then D22 text becomes shared string index 20
$ws.Range("D22").Value = "Added breeze and stench feedback into the learning algorith and ran some experiments to guage impact. Have to reward exploration otherwise Adventurer dies of starvation."

$ws.Range("A22").Value = 42917
$ws.Range("B22").Formula = "=A22"
$ws.Range("C22").Value = 2

$ws.Range("A23").Value = 42918
$ws.Range("B23").Formula = "=A23"
$ws.Range("C23").Value = 1.5

# Copy formatting from row 21 (A21:D21) to new rows 22 and 23, after values are set
$ws.Range("A21:D21").Copy()
$ws.Range("A22:D22").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A23:D23").PasteSpecial(-4122) # xlPasteFormats
$ws.Rows.Item(22).RowHeight = 28.5
$ws.Rows.Item(23).RowHeight = 28.5

[void]$ws.Range("D23").Select()

